# "select two cards at a time"
#
# The paragraph "2, display the cards." is split into four runs
# ("2, display ", "12 ", "cards", ", and make duplicates to 24 cards
# total.") that all share the original run's formatting, and two new
# paragraphs ("3, shuffle the cards" and "4, ") are inserted right
# after it, before the following (originally empty) paragraph.

$d = $word.ActiveDocument

# Locate the target sentence without disturbing anything else in the
# document. Find.Execute (with no replacement) collapses/extends the
# Range to exactly the matched text (the paragraph mark is excluded).
# NOTE: $d.Content returns a brand-new Range each time it is accessed,
# so it must be captured into a variable before calling Find.Execute
# on it; otherwise the match position can't be read back afterwards.
$findRange = $d.Content
$found = $findRange.Find.Execute("2, display the cards.", $true, $false, $false, `
                                  $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the paragraph to edit."
}
$startPos = $findRange.Start
$endPos = $findRange.End

# Build the replacement: four runs with the exact same run formatting
# as the original single run (Arial 13.5pt, white shading).
$runsXml = @'
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">2, display </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">12 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>cards</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>, and make duplicates to 24 cards total.</w:t></w:r>
'@

$payload = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>__RUNS__</w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$payload = $payload.Replace("__RUNS__", $runsXml)

# Using InsertXML (rather than Range.Text / Find.Replace) on a freshly
# constructed Range keeps the new runs distinct instead of letting the
# engine coalesce same-formatted runs together, and since the target
# range excludes the paragraph mark, the enclosing <w:p> element (and
# its attributes) are left untouched. A *fresh* Range object (built
# via Document.Range) is required here - reusing the Range that
# Find.Execute mutated in place causes InsertXML to append instead of
# replace.
$fresh = $d.Range($startPos, $endPos)
$fresh.InsertXML($payload)

# Re-derive the (now multi-run) paragraph that was just edited, by
# counting paragraph marks before it, so we can append two brand-new
# paragraphs right after it.
$paraIndex = $d.Range(0, $startPos).Paragraphs.Count + 1
$editedPara = $d.Paragraphs.Item($paraIndex)

$editedPara.Range.InsertParagraphAfter()
$thirdPara = $d.Paragraphs.Item($paraIndex + 1)
$thirdPara.Range.InsertAfter("3, shuffle the cards")

$thirdPara.Range.InsertParagraphAfter()
$fourthPara = $d.Paragraphs.Item($paraIndex + 2)
$fourthPara.Range.InsertAfter("4, ")
